$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 5 new data rows (rows 9-13) below the existing two data
#    rows (7 and 8). This pushes the old "totals" row (was row 9) to
#    row 14 and the old footer row (was row 10) to row 15.
# ------------------------------------------------------------------
$ws.Rows("9:13").Insert()

# Copy the formatting (styles, number formats, borders, fonts) from
# the template data row (row 8) into the newly inserted rows so they
# look exactly like the other item rows.
$ws.Range("A8:Q8").Copy()
$ws.Range("A9:Q13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Recreate the merged cells for each of the 5 new rows (A:B, C:G,
# H:K, L:M, N:O) matching the pattern used by rows 7 and 8.
for ($r = 9; $r -le 13; $r++) {
    $ws.Range("A" + $r + ":B" + $r).Merge()
    $ws.Range("C" + $r + ":G" + $r).Merge()
    $ws.Range("H" + $r + ":K" + $r).Merge()
    $ws.Range("L" + $r + ":M" + $r).Merge()
    $ws.Range("N" + $r + ":O" + $r).Merge()
}

# ------------------------------------------------------------------
# 2. Fill in the item rows (A = row number, C = item name, H = stock
#    ratio, L = order-limit flag, N = price, P = sell price, Q =
#    number-of-transactions ratio).
# ------------------------------------------------------------------
$items = @(
    @{ Row = 7;  Num = 1; Name = "ASPIRIN PROTECT 100MG 30 GASTRO-RESISTANT TAB"; H = "4:2"; L = "1"; N = "78.00";  P = "78.0000"; Q = "1:0" },
    @{ Row = 8;  Num = 2; Name = "CONCOR PLUS 5/12.5MG 30 F.C. TABLETS";          H = "1:3"; L = "1"; N = "72.00";  P = "72.0000"; Q = "1:0" },
    @{ Row = 9;  Num = 3; Name = "DIAMONRECTA 5 MG 30 F.C. TAB.";                 H = "1:1"; L = "1"; N = "187.50"; P = "61.8750"; Q = "0:1" },
    @{ Row = 10; Num = 4; Name = "LIBRAX 30 SUGAR COATED TAB";                    H = "9:0"; L = "1"; N = "48.00";  P = "15.8400"; Q = "0:1" },
    @{ Row = 11; Num = 5; Name = "MELOQUIN 4% CREAM 20 GM";                       H = "0:0"; L = "1"; N = "60.00";  P = "60.0000"; Q = "1:0" },
    @{ Row = 12; Num = 6; Name = "OXITROPIL 1200 MG 60 TAB";                      H = "0:3"; L = "1"; N = "123.00"; P = "40.5900"; Q = "0:2" },
    @{ Row = 13; Num = 7; Name = "URIPAN X.R. 10 MG 30 TAB.";                     H = "0:1"; L = "1"; N = "87.00";  P = "87.0000"; Q = "1:0" }
)

foreach ($it in $items) {
    $r = $it.Row
    $ws.Range("A" + $r).Value = $it.Num
    $ws.Range("C" + $r).Value = $it.Name
    $ws.Range("H" + $r).Value = $it.H
    $ws.Range("L" + $r).Value = $it.L
    $ws.Range("N" + $r).Value = $it.N
    $ws.Range("P" + $r).Value = $it.P
    $ws.Range("Q" + $r).Value = $it.Q
}

# ------------------------------------------------------------------
# 3. Update the totals row (now row 14) with the new sum of the
#    "sell price" column, and the footer timestamp row (now row 15).
# ------------------------------------------------------------------
$ws.Range("P14").Value = 415.30500000000001
$ws.Range("A15").Value = "Saturday, 19 July, 2025 10:17 AM"
